# Update the "data" sheet: the CBDB text list was regenerated, replacing
# the previous 13-row anthology list with a new 10-row list (commit:
# "include 稿 as anthology, fix bug for :").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# New tts_sysno (col A) / c_title_chn (col B) pairs, in final row order.
# Column C (c_created_by) stays constant at 57941 for every row.
$sysnos = @(562361, 562361, 562365, 562375, 562375, 562375, 72460, 72460, 72460, 72460)
$titles = @("鄉黨圖考", "詩集: 五卷", "周易闡要: 三卷", "蕉巖遺稿", "蓮塘詩稿", "元燈心法", "道古堂初刻", "道古齋詩稿", "道古齋文稿", "水經注校正")

for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $sysnos[$i]
    $ws.Cells.Item($row, 2).Value = $titles[$i]
    $ws.Cells.Item($row, 3).Value = 57941
}

# The old list had 13 rows; the new one only needs 10, so drop the tail.
$ws.Rows("11:13").Delete()

# Re-fit column B now that the text content (and its width) changed;
# 14.1667 is the input that lands on a stored/displayed width of 15.
$ws.Columns.Item(2).ColumnWidth = 14.1667

# Printable area is now explicitly portrait.
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# The "data" sheet (not "schema") is now the active/selected tab, with
# E4 as the selected cell.
$ws.Activate()
$ws.Range("E4").Select()
